$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Computer Name value
$ws.Range("C8").Value = "R08-LEY4-048"

# Mark Status column with check marks for completed items (all rows except 15)
$ws.Range("D14").Value = "✓"
$ws.Range("D16").Value = "✓"
$ws.Range("D17").Value = "✓"
$ws.Range("D18").Value = "✓"
$ws.Range("D19").Value = "✓"
$ws.Range("D20").Value = "✓"
$ws.Range("D21").Value = "✓"
$ws.Range("D22").Value = "✓"

# Turn on word-wrap for the "Problems Encountered/ Action" column entries
$ws.Range("E14:E22").WrapText = $true

# Update the view: scroll position and active selection
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("I15").Select()
